$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$c = $t.Cell(5, 3)
$r = $c.Range
$r.End = $r.End - 1
$r.Text = "-7.1"

$c = $t.Cell(5, 5)
$r = $c.Range
$r.End = $r.End - 1
$r.Text = "-4.0"

$c = $t.Cell(5, 7)
$r = $c.Range
$r.End = $r.End - 1
$r.Text = "8.0"

$c = $t.Cell(5, 9)
$r = $c.Range
$r.End = $r.End - 1
$r.Text = "-82.7"

$c = $t.Cell(5, 11)
$r = $c.Range
$r.End = $r.End - 1
$r.Text = "-14.6"

$c = $t.Cell(5, 13)
$r = $c.Range
$r.End = $r.End - 1
$r.Text = "-67.4"

$c = $t.Cell(6, 3)
$r = $c.Range
$r.End = $r.End - 1
$r.Text = "-8.2"

$c = $t.Cell(6, 5)
$r = $c.Range
$r.End = $r.End - 1
$r.Text = "-3.3"

$c = $t.Cell(6, 7)
$r = $c.Range
$r.End = $r.End - 1
$r.Text = "-4.6"

$c = $t.Cell(6, 9)
$r = $c.Range
$r.End = $r.End - 1
$r.Text = "-38.9"

$c = $t.Cell(6, 11)
$r = $c.Range
$r.End = $r.End - 1
$r.Text = "-6.7"

$c = $t.Cell(6, 13)
$r = $c.Range
$r.End = $r.End - 1
$r.Text = "-30.0"

$c = $t.Cell(7, 3)
$r = $c.Range
$r.End = $r.End - 1
$r.Text = "-22.2"

$c = $t.Cell(7, 5)
$r = $c.Range
$r.End = $r.End - 1
$r.Text = "-15.8"

$c = $t.Cell(7, 7)
$r = $c.Range
$r.End = $r.End - 1
$r.Text = "-35.9"

$c = $t.Cell(7, 9)
$r = $c.Range
$r.End = $r.End - 1
$r.Text = "-11.4"

$c = $t.Cell(7, 13)
$r = $c.Range
$r.End = $r.End - 1
$r.Text = "2.7"

$c = $t.Cell(8, 3)
$r = $c.Range
$r.End = $r.End - 1
$r.Text = "-10.6"

$c = $t.Cell(8, 5)
$r = $c.Range
$r.End = $r.End - 1
$r.Text = "-3.2"

$c = $t.Cell(8, 9)
$r = $c.Range
$r.End = $r.End - 1
$r.Text = "-28.6"

$c = $t.Cell(8, 11)
$r = $c.Range
$r.End = $r.End - 1
$r.Text = "7.4"

$c = $t.Cell(8, 13)
$r = $c.Range
$r.End = $r.End - 1
$r.Text = "-34.4"

$c = $t.Cell(9, 3)
$r = $c.Range
$r.End = $r.End - 1
$r.Text = "-5.2"

$c = $t.Cell(9, 5)
$r = $c.Range
$r.End = $r.End - 1
$r.Text = "-6.3"

$c = $t.Cell(9, 7)
$r = $c.Range
$r.End = $r.End - 1
$r.Text = "3.6"

$c = $t.Cell(9, 9)
$r = $c.Range
$r.End = $r.End - 1
$r.Text = "-8.8"

$c = $t.Cell(9, 11)
$r = $c.Range
$r.End = $r.End - 1
$r.Text = "6.4"

$c = $t.Cell(9, 13)
$r = $c.Range
$r.End = $r.End - 1
$r.Text = "-8.8"

$c = $t.Cell(10, 3)
$r = $c.Range
$r.End = $r.End - 1
$r.Text = "-55.8"

$c = $t.Cell(11, 3)
$r = $c.Range
$r.End = $r.End - 1
$r.Text = "-2.7"

$c = $t.Cell(11, 7)
$r = $c.Range
$r.End = $r.End - 1
$r.Text = "9.9"

$c = $t.Cell(11, 9)
$r = $c.Range
$r.End = $r.End - 1
$r.Text = "-12.2"

$c = $t.Cell(11, 11)
$r = $c.Range
$r.End = $r.End - 1
$r.Text = "-13.7"

$c = $t.Cell(11, 13)
$r = $c.Range
$r.End = $r.End - 1
$r.Text = "-14.9"

$c = $t.Cell(12, 3)
$r = $c.Range
$r.End = $r.End - 1
$r.Text = "-6.7"

$c = $t.Cell(12, 5)
$r = $c.Range
$r.End = $r.End - 1
$r.Text = "-6.8"

$c = $t.Cell(12, 7)
$r = $c.Range
$r.End = $r.End - 1
$r.Text = "-5.2"

$c = $t.Cell(12, 9)
$r = $c.Range
$r.End = $r.End - 1
$r.Text = "-66.2"

$c = $t.Cell(12, 11)
$r = $c.Range
$r.End = $r.End - 1
$r.Text = "-9.5"

$c = $t.Cell(13, 3)
$r = $c.Range
$r.End = $r.End - 1
$r.Text = "-8.0"

$c = $t.Cell(13, 9)
$r = $c.Range
$r.End = $r.End - 1
$r.Text = "-7.1"

$c = $t.Cell(13, 11)
$r = $c.Range
$r.End = $r.End - 1
$r.Text = "22.6"

$c = $t.Cell(14, 5)
$r = $c.Range
$r.End = $r.End - 1
$r.Text = "-7.2"

$c = $t.Cell(14, 9)
$r = $c.Range
$r.End = $r.End - 1
$r.Text = "-27.9"

$c = $t.Cell(14, 13)
$r = $c.Range
$r.End = $r.End - 1
$r.Text = "-13.4"

$c = $t.Cell(15, 3)
$r = $c.Range
$r.End = $r.End - 1
$r.Text = "-19.4"

$c = $t.Cell(15, 11)
$r = $c.Range
$r.End = $r.End - 1
$r.Text = "-18.9"

$c = $t.Cell(15, 13)
$r = $c.Range
$r.End = $r.End - 1
$r.Text = "-41.6"

$c = $t.Cell(16, 7)
$r = $c.Range
$r.End = $r.End - 1
$r.Text = "4.5"

$c = $t.Cell(16, 9)
$r = $c.Range
$r.End = $r.End - 1
$r.Text = "8.7"

$c = $t.Cell(16, 11)
$r = $c.Range
$r.End = $r.End - 1
$r.Text = "-17.5"

$c = $t.Cell(16, 13)
$r = $c.Range
$r.End = $r.End - 1
$r.Text = "-5.1"

$c = $t.Cell(18, 3)
$r = $c.Range
$r.End = $r.End - 1
$r.Text = "11.3"

$c = $t.Cell(18, 5)
$r = $c.Range
$r.End = $r.End - 1
$r.Text = "4.7"

$c = $t.Cell(18, 7)
$r = $c.Range
$r.End = $r.End - 1
$r.Text = "5.9"

$c = $t.Cell(18, 11)
$r = $c.Range
$r.End = $r.End - 1
$r.Text = "30.1"

$c = $t.Cell(19, 3)
$r = $c.Range
$r.End = $r.End - 1
$r.Text = "12.1"

$c = $t.Cell(19, 5)
$r = $c.Range
$r.End = $r.End - 1
$r.Text = "12.2"

$c = $t.Cell(19, 7)
$r = $c.Range
$r.End = $r.End - 1
$r.Text = "-12.7"

$c = $t.Cell(19, 11)
$r = $c.Range
$r.End = $r.End - 1
$r.Text = "13.5"

$c = $t.Cell(20, 7)
$r = $c.Range
$r.End = $r.End - 1
$r.Text = "13.1"

$c = $t.Cell(20, 9)
$r = $c.Range
$r.End = $r.End - 1
$r.Text = "-20.2"

$c = $t.Cell(20, 11)
$r = $c.Range
$r.End = $r.End - 1
$r.Text = "14.7"

$c = $t.Cell(20, 13)
$r = $c.Range
$r.End = $r.End - 1
$r.Text = "6.3"

$c = $t.Cell(21, 5)
$r = $c.Range
$r.End = $r.End - 1
$r.Text = "-2.5"

$c = $t.Cell(21, 7)
$r = $c.Range
$r.End = $r.End - 1
$r.Text = "5.2"

$c = $t.Cell(21, 9)
$r = $c.Range
$r.End = $r.End - 1
$r.Text = "-51.8"

$c = $t.Cell(21, 11)
$r = $c.Range
$r.End = $r.End - 1
$r.Text = "-7.1"

$c = $t.Cell(22, 3)
$r = $c.Range
$r.End = $r.End - 1
$r.Text = "4.6"

$c = $t.Cell(22, 5)
$r = $c.Range
$r.End = $r.End - 1
$r.Text = "11.3"

$c = $t.Cell(22, 7)
$r = $c.Range
$r.End = $r.End - 1
$r.Text = "12.1"

$c = $t.Cell(22, 9)
$r = $c.Range
$r.End = $r.End - 1
$r.Text = "9.6"

$c = $t.Cell(22, 11)
$r = $c.Range
$r.End = $r.End - 1
$r.Text = "-33.3"

$c = $t.Cell(22, 13)
$r = $c.Range
$r.End = $r.End - 1
$r.Text = "-19.2"

$c = $t.Cell(23, 7)
$r = $c.Range
$r.End = $r.End - 1
$r.Text = "-3.2"

$c = $t.Cell(23, 9)
$r = $c.Range
$r.End = $r.End - 1
$r.Text = "-28.7"

$c = $t.Cell(23, 11)
$r = $c.Range
$r.End = $r.End - 1
$r.Text = "-4.8"

$c = $t.Cell(23, 13)
$r = $c.Range
$r.End = $r.End - 1
$r.Text = "37.5"
